# Apply odds updates to Sheet1 for the 2026-02-04 Betfair Back/Lay workbook.
# Each line updates a single cell value per the upstream diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(3, 8).Value = 2.38
$ws.Cells.Item(3, 9).Value = 2.78
$ws.Cells.Item(3, 10).Value = 1.09
$ws.Cells.Item(3, 16).Value = 1.83
$ws.Cells.Item(3, 17).Value = 1.77
$ws.Cells.Item(4, 9).Value = 2.8
$ws.Cells.Item(4, 17).Value = 2.54
$ws.Cells.Item(4, 33).Value = 19
$ws.Cells.Item(8, 6).Value = 5.9
$ws.Cells.Item(8, 10).Value = 3.45
$ws.Cells.Item(8, 17).Value = 2.26
$ws.Cells.Item(9, 6).Value = 3.6
$ws.Cells.Item(9, 7).Value = 4.1
$ws.Cells.Item(9, 9).Value = 2.22
$ws.Cells.Item(9, 16).Value = 1.95
$ws.Cells.Item(10, 6).Value = 3.4
$ws.Cells.Item(10, 7).Value = 3.5
$ws.Cells.Item(10, 8).Value = 2.24
$ws.Cells.Item(10, 9).Value = 2.28
$ws.Cells.Item(10, 10).Value = 3.65
$ws.Cells.Item(10, 16).Value = 1.86
$ws.Cells.Item(10, 20).Value = 1.8
$ws.Cells.Item(10, 36).Value = 1000
$ws.Cells.Item(11, 6).Value = 4.5
$ws.Cells.Item(11, 7).Value = 4.8
$ws.Cells.Item(11, 8).Value = 1.92
$ws.Cells.Item(11, 9).Value = 1.99
$ws.Cells.Item(11, 10).Value = 3.5
$ws.Cells.Item(11, 17).Value = 1.91
$ws.Cells.Item(11, 24).Value = 970
$ws.Cells.Item(11, 26).Value = 970
$ws.Cells.Item(11, 27).Value = 23
$ws.Cells.Item(11, 29).Value = 970
$ws.Cells.Item(11, 31).Value = 24
$ws.Cells.Item(11, 32).Value = 34
$ws.Cells.Item(11, 37).Value = 70
$ws.Cells.Item(12, 9).Value = 13.5
$ws.Cells.Item(12, 14).Value = 5.2
$ws.Cells.Item(12, 16).Value = 2.44
$ws.Cells.Item(12, 19).Value = 2.5
$ws.Cells.Item(12, 21).Value = 1.81
$ws.Cells.Item(12, 24).Value = 26
$ws.Cells.Item(12, 27).Value = 550
$ws.Cells.Item(12, 28).Value = 9.6
$ws.Cells.Item(12, 32).Value = 8.6
$ws.Cells.Item(12, 33).Value = 11
$ws.Cells.Item(12, 34).Value = 32
$ws.Cells.Item(12, 36).Value = 10.5
$ws.Cells.Item(12, 40).Value = 4.8
$ws.Cells.Item(13, 8).Value = 5.5
$ws.Cells.Item(13, 16).Value = 2.06
$ws.Cells.Item(13, 21).Value = 2.06
$ws.Cells.Item(13, 24).Value = 19.5
$ws.Cells.Item(13, 27).Value = 170
$ws.Cells.Item(13, 28).Value = 10.5
$ws.Cells.Item(13, 29).Value = 9.4
$ws.Cells.Item(13, 30).Value = 25
$ws.Cells.Item(13, 35).Value = 85
$ws.Cells.Item(13, 36).Value = 21
$ws.Cells.Item(13, 37).Value = 21
$ws.Cells.Item(13, 39).Value = 120
$ws.Cells.Item(13, 40).Value = 1000
$ws.Cells.Item(15, 7).Value = 6.8
$ws.Cells.Item(15, 8).Value = 1.55
$ws.Cells.Item(15, 10).Value = 4.5
$ws.Cells.Item(15, 11).Value = 5
$ws.Cells.Item(15, 13).Value = 1.04
$ws.Cells.Item(15, 14).Value = 4.8
$ws.Cells.Item(15, 15).Value = 1.22
$ws.Cells.Item(15, 16).Value = 2.34
$ws.Cells.Item(15, 17).Value = 1.64
$ws.Cells.Item(15, 18).Value = 1.54
$ws.Cells.Item(15, 19).Value = 2.6
$ws.Cells.Item(15, 20).Value = 1.78
$ws.Cells.Item(15, 21).Value = 2.14
$ws.Cells.Item(15, 24).Value = 1000
$ws.Cells.Item(15, 25).Value = 10.5
$ws.Cells.Item(15, 29).Value = 12.5
$ws.Cells.Item(15, 31).Value = 18
$ws.Cells.Item(15, 32).Value = 1000
$ws.Cells.Item(15, 34).Value = 1000
$ws.Cells.Item(15, 35).Value = 1000
$ws.Cells.Item(15, 36).Value = 210
$ws.Cells.Item(16, 6).Value = 2.26
$ws.Cells.Item(16, 7).Value = 2.38
$ws.Cells.Item(16, 8).Value = 3.5
$ws.Cells.Item(16, 9).Value = 3.8
$ws.Cells.Item(17, 8).Value = 8.6
$ws.Cells.Item(17, 10).Value = 4.7
$ws.Cells.Item(17, 16).Value = 2.02
$ws.Cells.Item(19, 10).Value = 3.2
$ws.Cells.Item(20, 7).Value = 2.2
$ws.Cells.Item(20, 9).Value = 4.2
$ws.Cells.Item(20, 10).Value = 3.3
$ws.Cells.Item(20, 16).Value = 1.71
